# Apply updated betting odds values to row 6 of the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G6"  = 1.37
    "I6"  = 6.6
    "J6"  = 1.82
    "K6"  = 2.4
    "L6"  = 6
    "R6"  = 2.1
    "U6"  = 1.88
    "V6"  = 1.88
    "W6"  = 6.5
    "X6"  = 5.9
    "Z6"  = 7.7
    "AB6" = 19.5
    "AC6" = 13
    "AE6" = 15
    "AF6" = 60
    "AI6" = 35
    "AJ6" = 17
    "AK6" = 110
    "AL6" = 55
    "AN6" = 3.25
    "AO6" = 6.1
    "AP6" = 15
    "AQ6" = 16.5
    "AS6" = 175
    "AT6" = 3.2
    "AV6" = 65
    "AW6" = 8
    "AX6" = 37
    "AY6" = 37
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
